# Replace the literal "<br/>" marker with an actual line break in the
# four MSME size-class definitions (each value is shared between the
# "Size classification by annual turnover" tables in two places on the
# sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @{ Cells = @("C21", "C42"); Text = "<=500 SMMLV `n(Salario Minimo Mensual Legal Vigente - Legal monthly minimum wage in force)" },
    @{ Cells = @("C22", "C43"); Text = ">500 - <=5,000 SMMLV `n(Salario Minimo Mensual  Legal Vigente - Legal monthly minimum wage in force)" },
    @{ Cells = @("C23", "C44"); Text = ">5,000 - <=30,000 SMMLV `n(Salario Minimo Mensual Legal  Vigente - Legal monthly minimum wage in force)" },
    @{ Cells = @("C24", "C45"); Text = ">30,000 SMMLV `n(Salario Minimo Mensual Legal  Vigente - Legal monthly minimum wage in force)" }
)

foreach ($pair in $pairs) {
    foreach ($addr in $pair.Cells) {
        $ws.Range($addr).Value2 = $pair.Text
    }
}
